$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.562.11'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '2.900.61'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '528.08'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').Value = '143.51'
$ws.Range('E6').Value = '  -5.75%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Value = '2.909.95'
$ws.Range('E9').Value = '  -2.11%  '
$ws.Range('E10').Value = '  -4.17%  '
$ws.Range('D11').Value = '6.03'
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('D13').Value = '3.404.94'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '60.561.06'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '22.80'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -3.68%  '
$ws.Range('D17').Value = '2.904.68'
$ws.Range('E17').Value = '  -2.32%  '
$ws.Range('E18').Value = '  -3.91%  '
$ws.Range('D19').Value = '5.03'
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.70'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').Value = '361.73'
$ws.Range('E21').Value = '  -5.03%  '
$ws.Range('D22').Value = '6.64'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('E24').Value = '  -1.78%  '
$ws.Range('D25').Value = '64.83'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('D26').Value = '0.454'
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.180'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('D28').Value = '0.994'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').Value = '7.87'
$ws.Range('E29').Value = '  -5.38%  '
$ws.Range('D30').Value = '0.0₃0851'
$ws.Range('E30').Value = '  -8.58%  '
$ws.Range('E32').Value = '  -2.04%  '
$ws.Range('D33').Value = '19.74'
$ws.Range('E33').Value = '  -3.10%  '
$ws.Range('D34').Value = '152.12'
$ws.Range('E34').Value = '  -4.24%  '
$ws.Range('D35').Value = '4.38'
$ws.Range('E35').Value = '  -5.45%  '
$ws.Range('D36').Value = '5.58'
$ws.Range('E36').Value = '  -6.01%  '
$ws.Range('E37').Value = '  -5.58%  '
$ws.Range('E38').Value = '  -5.18%  '
$ws.Range('D39').Value = '37.65'
$ws.Range('E39').Value = '  +1.51%  '
$ws.Range('D40').Value = '1.48'
$ws.Range('E40').Value = '  -3.95%  '
$ws.Range('D41').Value = '3.72'
$ws.Range('E41').Value = '  -5.15%  '
$ws.Range('D42').Value = '2.293.84'
$ws.Range('E42').Value = '  -4.86%  '
$ws.Range('D43').Value = '0.648'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('D44').Value = '0.0581'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('D45').Value = '20.48'
$ws.Range('E45').Value = '  -7.91%  '
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').Value = '0.0238'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('E49').Value = '  -1.35%  '
$ws.Range('D50').Value = '0.0924'
$ws.Range('E50').Value = '  -3.13%  '
$ws.Range('D51').Value = '250.98'
$ws.Range('E51').Value = '  -6.31%  '
